# Nalco aluminium-ingot price sheet: daily refresh.
# A brand new "today" row is published (30-11-2025, same basic price/circular
# as the existing top row), so every historical row shifts down by one and
# the oldest row that falls off the bottom of the *old* range reappears as
# the new last row (117) purely because it was shifted there too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at the top of the data (row 2, just below the header).
# This automatically shifts every existing data row (2..116) down to (3..117)
# with their values, styles and formulas intact.
$ws.Rows.Item(2).Insert()

# The newly inserted row inherits a blank style from Insert(); repair it by
# pasting the formatting of the row right below (which holds the old row-2
# data) over it.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new top row. Use Value2 (not Value) so plain numbers/text go in
# without Excel's formula-bar autocorrect; date-like text still gets
# auto-parsed into a real date by Value2 when the day number is <=12 (mirrors
# normal Excel typed-entry behaviour for ambiguous D-M-Y strings), so those
# two cells are entered with a leading apostrophe to force literal text and
# then re-stamped with the neighbouring cell's format to drop the resulting
# "quote prefix" style variant.
$ws.Range("A2").Value2 = "'30-11-2025"
$ws.Range("B2").Value2 = "ALUMINIUM INGOT"
$ws.Range("C2").Value2 = "IE07"
$ws.Range("D2").Value2 = 297.15
$ws.Range("E2").Value2 = "'01-11-2025"
$ws.Range("F2").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

# Re-apply formats from row 3 to A2/E2 so they end up with the same style
# index as every other text cell (clears the quote-prefix style variant that
# the apostrophe-prefixed literal picked up).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# The new bottom row (117, ex-row 116 shifted down) needs its own hyperlink
# on column F, matching the other rows.
$ws.Hyperlinks.Add($ws.Range("F117"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Hyperlinks.Add stamps the cell with Excel's blue/underlined "Hyperlink"
# style; restore the plain data-row style by re-pasting formats from F116
# (identical row, just above) over it.
$ws.Range("F116").Copy()
$ws.Range("F117").PasteSpecial(-4122)
